$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Font.Bold = $false
$ws.Range("C1").Value = "average"
$ws.Range("C2").Value = "max"
$ws.Range("C3").Value = "min"

$ws.Range("D1").Formula = "=AVERAGE(A:A)"
$ws.Range("D2").Formula = "=MAX(A:A)"
$ws.Range("D3").Formula = "=MIN(A:A)"

$ws.Range("D1").Font.Bold = $true

$ws.Range("E13").Select()
